{"js": "// 1. Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateSearch = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateSearch.load(\"items\");\nawait context.sync();\nif (dateSearch.items.length > 0) {\n  dateSearch.items[0].insertText(\"September 21, 2025\", \"Replace\");\n  await context.sync();\n}\n\n// 2. Split the mailing-address paragraph \"999 Story Road, San Jose CA 95122\"\n//    into two paragraphs: \"999 Story Road\" and \"San Jose, CA 95122\" (outside\n//    of the property-address table, i.e. search only the document body's\n//    top-level paragraphs, not the table).\nconst addrSearch = context.document.body.search(\"999 Story Road, San Jose CA 95122\", { matchCase: true });\naddrSearch.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < addrSearch.items.length; i++) {\n  const range = addrSearch.items[i];\n  const para = range.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  // Only touch the standalone address paragraph, not the one that lives\n  // inside the property-address table.\n  const cell = para.parentTableCellOrNullObject;\n  cell.load(\"isNullObject\");\n  await context.sync();\n\n  if (!cell.isNullObject) {\n    continue;\n  }\n\n  // Replace the run's text with just the street line.\n  range.insertText(\"999 Story Road\", \"Replace\");\n  await context.sync();\n\n  // Insert a new paragraph after it carrying the city/state/zip line; it\n  // inherits the same paragraph/run formatting as the original paragraph.\n  para.insertParagraph(\"San Jose, CA 95122\", \"After\");\n  await context.sync();\n}\n\n// 3. Remove the empty \"NoSpacing\" paragraph directly under \"Board of Directors\".\nconst bodSearch = context.document.body.search(\"Board of Directors\", { matchCase: true });\nbodSearch.load(\"items\");\nawait context.sync();\n\nif (bodSearch.items.length > 0) {\n  const bodPara = bodSearch.items[0].paragraphs.getFirst();\n  const nextPara = bodPara.getNext();\n  nextPara.load(\"text,style\");\n  await context.sync();\n\n  if (nextPara.text.trim() === \"\") {\n    nextPara.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"September 19, 2025\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"September 21, 2025\"\n$find.Execute(\n    \"September 19, 2025\", $true, $false, $false, $false, $false,\n    $true, 1, $false, \"September 21, 2025\", 2\n) | Out-Null\n\n# 2. Split the mailing-address paragraph \"999 Story Road, San Jose CA 95122\"\n#    (the stand-alone paragraph in the address block -- NOT the copy that\n#    lives inside the property-address table) into two paragraphs:\n#    \"999 Story Road\" and \"San Jose, CA 95122\".\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    if ($r.Text.TrimEnd([char]13, [char]7) -eq \"999 Story Road, San Jose CA 95122\" -and $r.Information(12) -eq $false) {\n        # Information(12) == wdWithInTable; skip the table's copy of the address.\n        $target = $d.Range($r.Start, $r.End - 1)\n        $target.Text = \"999 Story Road\"\n        $target.InsertParagraphAfter()\n\n        $nextPara = $p.Next()\n        $nextRange = $d.Range($nextPara.Range.Start, $nextPara.Range.End - 1)\n        $nextRange.Text = \"San Jose, CA 95122\"\n        break\n    }\n}\n\n# 3. Remove the empty \"NoSpacing\" paragraph directly under \"Board of Directors\".\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Vietnam Town Condominium Owners Association Board of Directors\") {\n        $nextPara = $p.Next()\n        if ($nextPara.Range.Text.TrimEnd([char]13, [char]7) -eq \"\") {\n            $nextPara.Range.Delete() | Out-Null\n        }\n        break\n    }\n}\n"}
